# Excel file changes for Reorder Testcase
# Populates the DemoWebshop_ReOrder sheet (sheet4) with TestCase Description,
# UserName (with mailto hyperlink) and Password (with mailto hyperlink) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DemoWebshop_ReOrder")

$description = " DemoWebshop Application ReOrder"
$username = "aarosagarch@gmail.com"
$password = "Admin@123"

# Fill in the new columns C (Description), D (UserName), E (Password) for rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = $description
    $ws.Cells.Item($r, 4).Value = $username
    $ws.Cells.Item($r, 5).Value = $password
}

# Hyperlink the UserName column: D2 individually, D3:D11 as one linked range
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:aarosagarch@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3:D11"), "mailto:aarosagarch@gmail.com", "", "", "aarosagarch@gmail.com") | Out-Null

# Hyperlink the Password column: one hyperlink per cell, E2 through E11
for ($r = 2; $r -le 11; $r++) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), "mailto:Admin@123") | Out-Null
}

# Apply the built-in Hyperlink style to the UserName and Password columns
$ws.Range("D2:D11").Style = "Hyperlink"
$ws.Range("E2:E11").Style = "Hyperlink"

# Resize columns C and D to fit the new content (values chosen so the saved
# OOXML column width attribute lands on 34 and 23 respectively)
$ws.Columns(3).ColumnWidth = 33.166666666666664
$ws.Columns(4).ColumnWidth = 22.166666666666668

# Move the active selection to F11 on this sheet
$ws.Activate()
$ws.Range("F11").Select()
